$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-5 (34 columns, A:AH) with the new dataset ---
$data = New-Object 'object[,]' 4,34
$data[0,0] = 45155.50694444445
$data[0,1] = 23.06
$data[0,2] = 15.922
$data[0,3] = 4.233
$data[0,4] = 48.559
$data[0,5] = 40.17
$data[0,6] = 18.148
$data[0,7] = 60.258
$data[0,8] = 27.923
$data[0,9] = 11.884
$data[0,10] = 18.312
$data[0,11] = 19.215
$data[0,12] = 20.154
$data[0,13] = 5.794
$data[0,14] = 18.046
$data[0,15] = 25.401
$data[0,16] = 15.061
$data[0,17] = 3.797
$data[0,18] = 2.462
$data[0,19] = 267.342
$data[0,20] = 50.257
$data[0,21] = 16.657
$data[0,22] = 33.391
$data[0,23] = 17.393
$data[0,24] = 2.2
$data[0,25] = 29.943
$data[0,26] = 14.713
$data[0,27] = 13.198
$data[0,28] = 15.441
$data[0,29] = 19.895
$data[0,30] = 3.64
$data[0,31] = 53.246
$data[0,32] = 9.286
$data[0,33] = 20.825
$data[1,0] = 45155.51388888889
$data[1,1] = 8.648
$data[1,2] = 5.814
$data[1,3] = 1.556
$data[1,4] = 18.177
$data[1,5] = 15.061
$data[1,6] = 6.806
$data[1,7] = 30.22
$data[1,8] = 10.471
$data[1,9] = 4.408
$data[1,10] = 6.728
$data[1,11] = 7.299
$data[1,12] = 7.531
$data[1,13] = 2.178
$data[1,14] = 6.767
$data[1,15] = 9.487
$data[1,16] = 5.909
$data[1,17] = 1.563
$data[1,18] = 0.827
$data[1,19] = 95.70699999999999
$data[1,20] = 19.109
$data[1,21] = 6.247
$data[1,22] = 12.514
$data[1,23] = 6.73
$data[1,24] = 0.755
$data[1,25] = 14.004
$data[1,26] = 5.517
$data[1,27] = 5.1
$data[1,28] = 5.947
$data[1,29] = 7.548
$data[1,30] = 1.294
$data[1,31] = 27.439
$data[1,32] = 3.401
$data[1,33] = 7.81
$data[2,0] = 45155.52083333334
$data[2,1] = 23.541
$data[2,2] = 17.235
$data[2,3] = 1.562
$data[2,4] = 50.797
$data[2,5] = 42.026
$data[2,6] = 18.526
$data[2,7] = 69.193
$data[2,8] = 28.505
$data[2,9] = 12.591
$data[2,10] = 18.895
$data[2,11] = 20.447
$data[2,12] = 21.449
$data[2,13] = 5.917
$data[2,14] = 18.422
$data[2,15] = 26.173
$data[2,16] = 15.528
$data[2,17] = 1.146
$data[2,18] = 1.019
$data[2,19] = 273.086
$data[2,20] = 51.456
$data[2,21] = 17.004
$data[2,22] = 34.537
$data[2,23] = 18.36
$data[2,24] = 2.305
$data[2,25] = 33.881
$data[2,26] = 15.02
$data[2,27] = 13.365
$data[2,28] = 15.695
$data[2,29] = 21.411
$data[2,30] = 0.784
$data[2,31] = 62.472
$data[2,32] = 9.568
$data[2,33] = 21.259
$data[3,0] = 45155.52777777778
$data[3,1] = 2.88
$data[3,2] = 1.83
$data[3,3] = 0.65
$data[3,4] = 5.96
$data[3,5] = 4.95
$data[3,6] = 2.27
$data[3,7] = 16.64
$data[3,8] = 3.49
$data[3,9] = 1.52
$data[3,10] = 2.15
$data[3,11] = 2.47
$data[3,12] = 2.44
$data[3,13] = 0.75
$data[3,14] = 2.26
$data[3,15] = 3.26
$data[3,16] = 2.1
$data[3,17] = 0.72
$data[3,18] = 0.31
$data[3,19] = 27.12
$data[3,20] = 6.8
$data[3,21] = 2.08
$data[3,22] = 4.47
$data[3,23] = 2.41
$data[3,24] = 0.2
$data[3,25] = 7.28
$data[3,26] = 1.84
$data[3,27] = 1.77
$data[3,28] = 2.05
$data[3,29] = 2.53
$data[3,30] = 0.5600000000000001
$data[3,31] = 15.57
$data[3,32] = 1.08
$data[3,33] = 2.61

$ws.Range("A2:AH5").Value = $data

# --- Remove old row 6 (dataset now has one fewer row) ---
$ws.Rows.Item(6).Delete()

# --- Column width tweaks ---
$ws.Columns.Item(3).ColumnWidth = 8 - 5/6
$ws.Columns.Item(10).ColumnWidth = 8 - 5/6
$ws.Columns.Item(11).ColumnWidth = 8 - 5/6
$ws.Columns.Item(17).ColumnWidth = 8 - 5/6
$ws.Columns.Item(20).ColumnWidth = 9 - 5/6
$ws.Columns.Item(21).ColumnWidth = 8 - 5/6
$ws.Columns.Item(27).ColumnWidth = 8 - 5/6
$ws.Columns.Item(28).ColumnWidth = 8 - 5/6
$ws.Columns.Item(29).ColumnWidth = 8 - 5/6
